# commit: "fixed bugs #1 - #4"
#
# The test log sheet gets a new "EdgeTester::testEdge01" row (mirroring the
# existing Test Suite / Test Case / Observed Failure / Fix rows already on
# row 3 & row 4) plus a new "#4 ..." fix note appended under the "Fixes:"
# section at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: new Test Suite / Test Case / Observed Failure / Fix entry for bug #4
$ws.Range("B6").Value = "EdgeTester::testEdge01"
$ws.Range("C6").Value = "testing e.getSlopeZ"
$ws.Range("D6").Value = "not setting value to inf"
$ws.Range("E6").Value = "#4"

# New fix note appended below the existing #1/#2/#3 entries
$ws.Range("B26").Value = '#4 changed "!=" to "==" Edge.cpp line 94'

# Scroll/select so the newly-added row is in view, same as the author
# leaving the cursor on the last edited cell before saving.
$ws.Range("B26").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
